$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the timestamp value on existing row 6
$ws.Cells.Item(6, 1).Value = 45729.51109761574

$recommendation = "`nRecommendations for calcium Deficiency:`nTofu, raw, firm, prepared with calcium sulfate`nCheese, Mexican, blend, reduced fat`nCheese, cheddar, nonfat or fat free`nCheese, swiss, low fat`nCheese, swiss, low sodium`nCheese, mozzarella, part skim milk`nCheese, gruyere`nCheese, monterey`nCheese, port de salut`nCheese, swiss`nCheese, swiss`nCheese, provolone, sliced`nCheese, provolone, reduced fat`nCheese, monterey jack, solid`nCheese, low-sodium, cheddar or colby`nCheese, muenster`nCheese, mozzarella, low sodium`nCheese, provolone`nCheese, monterey, low fat`nCheese, brick`nCheese, mexican, queso asadero`nCheese, colby`nCheese, Mexican blend`nCheese, Swiss, nonfat or fat free`nCheese, queso fresco, solid`nCheese, cheddar`nCheese, mexican, queso chihuahua`nCheese, cheddar, sharp, sliced`nCheese, cheddar`nCheese, white, queso blanco`nCheese, mozzarella, nonfat`nCheese, cheddar, reduced fat`nCheese, tilsit`nCheese, parmesan, grated, refrigerated`nCheese, cheshire`nCheese, parmesan, hard`nCheese, caraway`nImitation cheese, american or cheddar, low cholesterol`nCheese, fontina`nCheese, mexican, queso anejo"

$newRows = @(
    @{ Row = 7; Timestamp = 45730.43464201389 },
    @{ Row = 8; Timestamp = 45730.43493046296 },
    @{ Row = 9; Timestamp = 45730.43550134636 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item(6, 1).NumberFormat
    $ws.Cells.Item($row, 1).Value = $r.Timestamp
    $ws.Cells.Item($row, 2).Value = "Aarti"
    $ws.Cells.Item($row, 3).Value = 25
    $ws.Cells.Item($row, 4).Value = "Female"
    $ws.Cells.Item($row, 5).Value = 50
    $ws.Cells.Item($row, 6).Value = 1.5
    $ws.Cells.Item($row, 7).Value = 22.22
    $ws.Cells.Item($row, 8).Value = "Normal weight - Maintain a balanced diet and exercise."
    $ws.Cells.Item($row, 9).Value = "Veg"
    $ws.Cells.Item($row, 10).Value = "calcium"
    $ws.Cells.Item($row, 11).Value = $recommendation
    $ws.Rows.Item($row).AutoFit()
}
